# Apply updated leve-profit figures (currentAveragePrice / profit columns)
# recomputed by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2447.4487
$ws.Range("I2").Value = 1043.5555
$ws.Range("K2").Value = 1043.5555
$ws.Range("M2").Value = -930.5554999999999
# Row 34
$ws.Range("H34").Value = 5024999.5
$ws.Range("I34").Value = 5024999.5
$ws.Range("K34").Value = 5024999.5
$ws.Range("M34").Value = -5024728.5
# Row 61
$ws.Range("H61").Value = 5638.4097
$ws.Range("I61").Value = 3702.0344
$ws.Range("K61").Value = 3702.0344
$ws.Range("M61").Value = -3490.0344
# Row 110
$ws.Range("H110").Value = 5897.967
$ws.Range("I110").Value = 2813.889
$ws.Range("K110").Value = 2813.889
$ws.Range("M110").Value = -768.8890000000001
# Row 116
$ws.Range("H116").Value = 2447.4487
$ws.Range("I116").Value = 1043.5555
$ws.Range("K116").Value = 1043.5555
$ws.Range("M116").Value = 1250.4445
# Row 132
$ws.Range("H132").Value = 527357.8
$ws.Range("J132").Value = 96193.7
$ws.Range("L132").Value = 288581.1
$ws.Range("N132").Value = -293641.1
# Row 136
$ws.Range("H136").Value = 5638.4097
$ws.Range("I136").Value = 3702.0344
$ws.Range("K136").Value = 11106.1032
$ws.Range("M136").Value = -8556.1032

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2447.4487
$ws.Range("I3").Value = 1043.5555
$ws.Range("K3").Value = 1043.5555
$ws.Range("M3").Value = -929.5554999999999
# Row 26
$ws.Range("H26").Value = 5592.5
$ws.Range("I26").Value = 5592.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5592.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -5300.5
$ws.Range("N26").Value = $null
# Row 99
$ws.Range("H99").Value = 8832.940000000001
$ws.Range("I99").Value = 9631.5
$ws.Range("J99").Value = 8010.1816
$ws.Range("K99").Value = 9631.5
$ws.Range("L99").Value = 8010.1816
$ws.Range("M99").Value = -8133.5
$ws.Range("N99").Value = -11006.1816
# Row 105
$ws.Range("H105").Value = 2192.56
$ws.Range("I105").Value = 1974.4783
$ws.Range("K105").Value = 1974.4783
$ws.Range("M105").Value = -227.4783
# Row 134
$ws.Range("H134").Value = 870510.0600000001
$ws.Range("J134").Value = 9193.700000000001
$ws.Range("L134").Value = 27581.1
$ws.Range("N134").Value = -32651.1

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 7205.25
$ws.Range("I58").Value = 4432
$ws.Range("J58").Value = 13059.889
$ws.Range("K58").Value = 4432
$ws.Range("L58").Value = 13059.889
$ws.Range("M58").Value = -4229
$ws.Range("N58").Value = -13465.889
# Row 134
$ws.Range("H134").Value = 11436.723
$ws.Range("I134").Value = 7450.875
$ws.Range("J134").Value = 14625.4
$ws.Range("K134").Value = 22352.625
$ws.Range("L134").Value = 43876.2
$ws.Range("M134").Value = -19817.625
$ws.Range("N134").Value = -48946.2
# Row 136
$ws.Range("H136").Value = 7205.25
$ws.Range("I136").Value = 4432
$ws.Range("J136").Value = 13059.889
$ws.Range("K136").Value = 13296
$ws.Range("L136").Value = 39179.667
$ws.Range("M136").Value = -10746
$ws.Range("N136").Value = -44279.667

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 111935.55
$ws.Range("J37").Value = 111935.55
$ws.Range("L37").Value = 335806.65
$ws.Range("N37").Value = -336030.65
# Row 63
$ws.Range("H63").Value = 12762.2
$ws.Range("I63").Value = 1962
$ws.Range("J63").Value = 19962.334
$ws.Range("K63").Value = 5886
$ws.Range("L63").Value = 59887.00199999999
$ws.Range("M63").Value = -5137
$ws.Range("N63").Value = -61385.00199999999
# Row 66
$ws.Range("H66").Value = 12762.2
$ws.Range("I66").Value = 1962
$ws.Range("J66").Value = 19962.334
$ws.Range("K66").Value = 17658
$ws.Range("L66").Value = 179661.006
$ws.Range("M66").Value = -13914
$ws.Range("N66").Value = -187149.006
# Row 122
$ws.Range("H122").Value = 3559.9375
$ws.Range("I122").Value = 195
$ws.Range("J122").Value = 4040.6428
$ws.Range("K122").Value = 1755
$ws.Range("L122").Value = 36365.7852
$ws.Range("M122").Value = 695
$ws.Range("N122").Value = -41265.7852

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 7585.2666
$ws.Range("I113").Value = 3420
$ws.Range("J113").Value = 13833.167
$ws.Range("K113").Value = 3420
$ws.Range("L113").Value = 13833.167
$ws.Range("M113").Value = -1250
$ws.Range("N113").Value = -18173.167
# Row 122
$ws.Range("H122").Value = 7171.615
$ws.Range("I122").Value = 7089
$ws.Range("J122").Value = 7242.4287
$ws.Range("K122").Value = 21267
$ws.Range("L122").Value = 21727.2861
$ws.Range("M122").Value = -18817
$ws.Range("N122").Value = -26627.2861
# Row 126
$ws.Range("H126").Value = 55574376
$ws.Range("I126").Value = 500000000
$ws.Range("J126").Value = 21174.125
$ws.Range("K126").Value = 1500000000
$ws.Range("L126").Value = 63522.375
$ws.Range("M126").Value = -1499997530
$ws.Range("N126").Value = -68462.375
# Row 132
$ws.Range("H132").Value = 4399.745
$ws.Range("I132").Value = 4117.0264
$ws.Range("J132").Value = 5226.154
$ws.Range("K132").Value = 12351.0792
$ws.Range("L132").Value = 15678.462
$ws.Range("M132").Value = -9821.0792
$ws.Range("N132").Value = -20738.462

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 27779670
$ws.Range("J46").Value = 38463804
$ws.Range("L46").Value = 38463804
$ws.Range("N46").Value = -38464180
# Row 122
$ws.Range("H122").Value = 1338702.6
$ws.Range("I122").Value = 2003953.8
$ws.Range("J122").Value = 8200.4
$ws.Range("K122").Value = 6011861.4
$ws.Range("L122").Value = 24601.2
$ws.Range("M122").Value = -6009411.4
$ws.Range("N122").Value = -29501.2
# Row 132
$ws.Range("H132").Value = 6959.077
$ws.Range("I132").Value = 6049.6
$ws.Range("J132").Value = 8199.272000000001
$ws.Range("K132").Value = 18148.8
$ws.Range("L132").Value = 24597.816
$ws.Range("M132").Value = -15618.8
$ws.Range("N132").Value = -29657.816
# Row 140
$ws.Range("H140").Value = 59357
$ws.Range("J140").Value = 58416.5
$ws.Range("L140").Value = 58416.5
$ws.Range("N140").Value = -68776.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 6240.729
$ws.Range("I132").Value = 5973
$ws.Range("J132").Value = 6686.9443
$ws.Range("K132").Value = 17919
$ws.Range("L132").Value = 20060.8329
$ws.Range("M132").Value = -15389
$ws.Range("N132").Value = -25120.8329
# Row 136
$ws.Range("H136").Value = 7814.9663
$ws.Range("I136").Value = 7581.45
$ws.Range("K136").Value = 22744.35
$ws.Range("M136").Value = -20194.35
